# The source commit only touches the PowerPoll "OfficeApp" web-extension
# content app embedded on slide 1 (ppt/slides/udata/data.xml + the
# graphicFrame/picture fallback pair in ppt/slides/slide.xml) - the pie
# chart there is rendered live by that third-party JS task-pane add-in,
# and its colors are state that lives inside the add-in's own runtime,
# not in any property the PowerPoint object model exposes. There is no
# WebExtension/AddIn-content object, method, or property in the Shape,
# Shapes, Slide, or Presentation COM interfaces that reaches the
# we:webextension part or its snapshot image, so that part of the diff
# (the webextension id GUID and the embedded snapshot picture bytes)
# cannot be produced through COM automation - exactly like real
# PowerPoint, where only the add-in itself (via Office.js, not VBA/COM)
# can re-render and re-snapshot its chart.
#
# The rest of the diff (every r:id on sldMasterId/sldId/sldLayoutId and
# on the webextensionref/blip relationships) is just fresh random
# relationship ids stamped in by the tool that resaved the package after
# that edit - it touches parts that have nothing to do with the pie
# chart (slide master, slide layouts, the slide relationship itself) and
# carries no semantic content of its own.
#
# So this script walks the reachable slide/shape object model (the only
# placeholders actually exposed here are the title and subtitle, neither
# of which changed in the source commit) to confirm the deck loads and
# leaves it untouched, rather than emitting destructive placeholder
# operations (e.g. Shapes.AddPicture/Fill.UserPicture with no real image
# bytes available to this script) that would corrupt the one shape that
# *is* trying to represent the add-in's content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

Write-Output "Slides: $($p.Slides.Count)"
Write-Output "Shapes on slide 1: $($s.Shapes.Count)"
